$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 12 (CONTINUE_MAIN_TEST) before removing it, then delete the
# entire row, shifting everything below it up by one (rows 13-17 -> 12-16).
[void]$ws.Rows(12).Select()
$ws.Rows(12).Delete()

# Rename the remaining two dictionary keys that were renamed in this revision.
$ws.Range("A2").Value = "TESTNAME"
$ws.Range("A6").Value = "FEEDBACK"
